# Update the "Report" sheet data: column D (Mango) values are pending
# review and have been zeroed out; columns E-G adjusted; column B
# ("General") recalculated as the row sum of C:H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

$data = @{
    2  = @(2657, 704, 0, 45, 622, 543, 743)
    3  = @(1133, 343, 0, 36, 265, 216, 273)
    4  = @(1419, 258, 0, 64, 196, 120, 781)
    5  = @(654, 218, 0, 32, 152, 139, 113)
    6  = @(499, 181, 0, 17, 89, 43, 169)
    7  = @(282, 20, 0, 4, 150, 42, 66)
    8  = @(647, 146, 0, 295, 131, 58, 17)
    9  = @(207, 43, 0, 95, 56, 8, 5)
    10 = @(616, 115, 0, 22, 339, 80, 60)
    11 = @(185, 9, 0, 0, 0, 176, 0)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 0; $col -lt $values.Length; $col++) {
        # Columns B..H correspond to array indices 0..6 -> spreadsheet columns 2..8
        $ws.Cells.Item($row, $col + 2).Value = $values[$col]
    }
}
